$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Newslatter")

$ws.Range("A25").Value = "noeloo@o2.pl"
$ws.Range("A26").Value = "michaltak830@gmail.com"
$ws.Range("A27").Value = "takk@o2.pl"
